$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-04-09"

# Update the header label cell (I1) which holds the "2022 (through 04-08)" text
$ws.Range("I1").Value = "2022 (through 04-09)"

# Update the April (row 5) 2022-YTD figure
$ws.Range("I5").Value = 32

# Update the Total (row 14) 2022-YTD figure
$ws.Range("I14").Value = 466
